$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [Coin, Link, Price, Volume(1h)] for rows 2..51 of the sheet.
# All four columns (B:E) are stored as text in the source data, so we
# force NumberFormat = "@" before writing, otherwise values that look
# numeric (e.g. "245.10", "1.000", "31.00") would be coerced to numbers
# and lose their original formatting / trailing zeros.
$data = @{
    2 = @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.495.38', '  +0.56%  ')
    3 = @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.727.80', '  +0.60%  ')
    4 = @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9995', '  +0.10%  ')
    5 = @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '245.10', '  +2.39%  ')
    6 = @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9998', '  +0.03%  ')
    7 = @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4802', '  +1.96%  ')
    8 = @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2678', '  +1.48%  ')
    9 = @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06218', '  -0.08%  ')
    10 = @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.727.25', '  +0.68%  ')
    11 = @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07149', '  +1.06%  ')
    12 = @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.69', '  +2.98%  ')
    13 = @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6185', '  +4.82%  ')
    14 = @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.514', '  +2.12%  ')
    15 = @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.18', '  +1.08%  ')
    16 = @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9999', '  -0.01%  ')
    17 = @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.502.06', '  +0.67%  ')
    18 = @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.09%  ')
    19 = @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006939', '  +1.81%  ')
    20 = @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.67', '  +1.09%  ')
    21 = @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.949.83', '  +0.98%  ')
    22 = @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.535', '  -0.41%  ')
    23 = @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.959', '  +1.54%  ')
    24 = @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.294', '  -1.08%  ')
    25 = @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '136.51', '  +0.83%  ')
    26 = @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.34', '  +0.88%  ')
    27 = @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.801', '  +2.17%  ')
    28 = @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.408', '  +0.10%  ')
    29 = @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '106.80', '  -0.17%  ')
    30 = @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.975', '  -1.82%  ')
    31 = @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08013', '  +3.86%  ')
    32 = @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.720', '  +0.82%  ')
    33 = @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04561', '  +3.22%  ')
    34 = @('Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9994', '  +0.00%  ')
    35 = @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.615', '  +0.15%  ')
    36 = @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6368', '  +2.23%  ')
    37 = @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.9897', '  +1.81%  ')
    38 = @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9341', '  +0.91%  ')
    39 = @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.089', '  +9.55%  ')
    40 = @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.417', '  -0.22%  ')
    41 = @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '104.89', '  -7.50%  ')
    42 = @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.003', '  +0.26%  ')
    43 = @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01503', '  +2.40%  ')
    44 = @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.679', '  +7.83%  ')
    45 = @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3907', '  +2.52%  ')
    46 = @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.920', '  +10.84%  ')
    47 = @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1189', '  +3.40%  ')
    48 = @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05330', '  +0.79%  ')
    49 = @('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '31.00', '  +1.59%  ')
    50 = @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.870', '  +2.50%  ')
    51 = @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.268', '  +3.82%  ')
}

foreach ($r in $data.Keys) {
    $rowIdx = [int]$r
    $vals = $data[$r]
    for ($col = 2; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($rowIdx, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col - 2]
    }
}
